# Update "想去人数" (want-to-go count) values in column F for both the
# "展览" (Exhibition) sheet and the "全部类型" (All types) aggregate sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Row -> New value, for sheet "展览"
$sheet1Updates = @{
    4  = 1441
    5  = 182
    6  = 36
    8  = 9670
    9  = 163
    12 = 185
    14 = 6662
    15 = 1081
    16 = 123
    18 = 180
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# Row -> New value, for sheet "全部类型"
$sheet4Updates = @{
    4  = 1441
    5  = 182
    6  = 36
    10 = 9670
    11 = 163
    14 = 185
    16 = 6662
    17 = 1081
    18 = 123
    20 = 180
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
